# Rename the worksheet from "Property1" to "DataNode" to unify the
# DataNode / DataTable / Entity naming convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move the live selection/active cell the same way the author's Excel
# session left it after editing (sheet view state).
[void]$ws.Range("E23").Select()
